# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.952.83'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '3.741.12'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''600.65'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '''165.38'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('D7').Value = '3.738.16'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').Value = '''0.171'
$ws.Range('E10').Value = '  +4.19%  '
$ws.Range('D11').Value = '''6.40'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('D12').Value = '''0.458'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '''37.67'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = '4.365.69'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '3.738.57'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '69.023.17'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('D19').Value = '''17.67'
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '''11.28'
$ws.Range('E21').Value = '  +4.55%  '
$ws.Range('D22').Value = '''490.12'
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').Value = '''0.724'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').Value = '''84.53'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').Value = '''0.0000148'
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('E26').Value = '  -2.09%  '
$ws.Range('D27').Value = '''12.32'
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('D28').Value = '''10.05'
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').Value = '''8.15'
$ws.Range('E31').Value = '  +2.94%  '
$ws.Range('D32').Value = '''2.44'
$ws.Range('E32').Value = '  -5.48%  '
$ws.Range('D33').Value = '''31.54'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('D34').Value = '3.884.54'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '3.675.80'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('E39').Value = '  +3.75%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('D41').Value = '''3.12'
$ws.Range('E41').Value = '  +7.78%  '
$ws.Range('D42').Value = '''0.324'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').Value = '''48.61'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').Value = '''426.90'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '''1.31'
$ws.Range('E48').Value = '  +9.99%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = '''39.93'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''141.18'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = '2.785.59'
$ws.Range('E51').Value = '  -0.28%  '
